# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (and a few cell additions/removals) to the
# "Leve Profits" worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as
# described by the source diff for Sheets/Exodus_Profits.xlsx.

$wb = $excel.ActiveWorkbook

# ----- Worksheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3249.2856
$ws.Range("I70").Value = 3016
$ws.Range("J70").Value = 3424.25
$ws.Range("K70").Value = 9048
$ws.Range("L70").Value = 10272.75
$ws.Range("M70").Value = -8778
$ws.Range("N70").Value = -10812.75

$ws.Range("H73").Value = 3249.2856
$ws.Range("I73").Value = 3016
$ws.Range("J73").Value = 3424.25
$ws.Range("K73").Value = 9048
$ws.Range("L73").Value = 10272.75
$ws.Range("M73").Value = -8112
$ws.Range("N73").Value = -12144.75

$ws.Range("H96").Value = 91501.27
$ws.Range("I96").Value = 143416.58
$ws.Range("J96").Value = 649.5
$ws.Range("K96").Value = 430249.74
$ws.Range("L96").Value = 1948.5
$ws.Range("M96").Value = -428876.74
$ws.Range("N96").Value = -4694.5

$ws.Range("H129").Value = 2277.8572
$ws.Range("I129").Value = 1431.3334
$ws.Range("K129").Value = 4294.0002
$ws.Range("M129").Value = 705.9997999999996

$ws.Range("H132").Value = 1968.2452
$ws.Range("I132").Value = 1732.7441
$ws.Range("K132").Value = 5198.2323
$ws.Range("M132").Value = -2668.2323

$ws.Range("H137").Value = 727271.7
$ws.Range("I137").Value = 2574
$ws.Range("J137").Value = 1613013.4
$ws.Range("K137").Value = 7722
$ws.Range("L137").Value = 4839040.199999999
$ws.Range("M137").Value = -5172
$ws.Range("N137").Value = -4844140.199999999

$ws.Range("H138").Value = 1772.5532
$ws.Range("I138").Value = 1379.7894
$ws.Range("J138").Value = 2039.0714
$ws.Range("K138").Value = 4139.3682
$ws.Range("L138").Value = 6117.2142
$ws.Range("M138").Value = 1000.6318
$ws.Range("N138").Value = -16397.2142


# ----- Worksheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2663.5293
$ws.Range("I45").Value = 2118.6365
$ws.Range("K45").Value = 2118.6365
$ws.Range("M45").Value = -1741.6365

$ws.Range("H74").Value = 2776.2173
$ws.Range("I74").Value = 2047.375
$ws.Range("J74").Value = 3164.9333
$ws.Range("K74").Value = 2047.375
$ws.Range("L74").Value = 3164.9333
$ws.Range("M74").Value = -1173.375
$ws.Range("N74").Value = -4912.933300000001

$ws.Range("H77").Value = 2776.2173
$ws.Range("I77").Value = 2047.375
$ws.Range("J77").Value = 3164.9333
$ws.Range("K77").Value = 10236.875
$ws.Range("L77").Value = 15824.6665
$ws.Range("M77").Value = -5868.875
$ws.Range("N77").Value = -24560.6665

$ws.Range("H110").Value = 932.17645
$ws.Range("I110").Value = 865.4375
$ws.Range("K110").Value = 865.4375
$ws.Range("M110").Value = 1179.5625

$ws.Range("H132").Value = 2302.5715
$ws.Range("I132").Value = 1476.7222
$ws.Range("K132").Value = 4430.1666
$ws.Range("M132").Value = -1900.1666

$ws.Range("H135").Value = 96091.75
$ws.Range("J135").Value = 96091.75
$ws.Range("L135").Value = 96091.75
$ws.Range("N135").Value = -106231.75

$ws.Range("H140").Value = 87970.664
$ws.Range("J140").Value = 87970.664
$ws.Range("L140").Value = 87970.664
$ws.Range("N140").Value = -98330.664


# ----- Worksheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 38702.297
$ws.Range("J105").Value = 3069.2
$ws.Range("L105").Value = 3069.2
$ws.Range("N105").Value = -6563.2


# ----- Worksheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1957.35
$ws.Range("I16").Value = 1592.6923
$ws.Range("J16").Value = 2634.5715
$ws.Range("K16").Value = 1592.6923
$ws.Range("L16").Value = 2634.5715
$ws.Range("M16").Value = -1305.6923
$ws.Range("N16").Value = -3208.5715

$ws.Range("H31").Value = 2321.9092
$ws.Range("I31").Value = 1629.9615
$ws.Range("K31").Value = 1629.9615
$ws.Range("M31").Value = -1334.9615

$ws.Range("H34").Value = 2321.9092
$ws.Range("I34").Value = 1629.9615
$ws.Range("K34").Value = 1629.9615
$ws.Range("M34").Value = -1427.9615

$ws.Range("H99").Value = 3128973.2
$ws.Range("I99").Value = 4172.1665
$ws.Range("K99").Value = 4172.1665
$ws.Range("M99").Value = -2674.1665

$ws.Range("H113").Value = 1957.35
$ws.Range("I113").Value = 1592.6923
$ws.Range("J113").Value = 2634.5715
$ws.Range("K113").Value = 1592.6923
$ws.Range("L113").Value = 2634.5715
$ws.Range("M113").Value = 577.3077000000001
$ws.Range("N113").Value = -6974.5715

$ws.Range("H126").Value = 3128973.2
$ws.Range("I126").Value = 4172.1665
$ws.Range("K126").Value = 12516.4995
$ws.Range("M126").Value = -10046.4995

$ws.Range("H132").Value = 1651.4546
$ws.Range("I132").Value = 1158.25
$ws.Range("K132").Value = 3474.75
$ws.Range("M132").Value = -944.75

$ws.Range("H141").Value = 131529.67
$ws.Range("J141").Value = 134688
$ws.Range("L141").Value = 134688
$ws.Range("N141").Value = -145048


# ----- Worksheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3729.75
$ws.Range("I46").Value = 4543.1665
$ws.Range("J46").Value = 1289.5
$ws.Range("K46").Value = 13629.4995
$ws.Range("L46").Value = 3868.5
$ws.Range("M46").Value = -13538.4995
$ws.Range("N46").Value = -4050.5

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H94").Value = 2079.8
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 2949.5
$ws.Range("K94").Value = 4500
$ws.Range("L94").Value = 8848.5
$ws.Range("M94").Value = -3824
$ws.Range("N94").Value = -10200.5

$ws.Range("H116").Value = 1380.5
$ws.Range("I116").Value = 1398.25
$ws.Range("K116").Value = 4194.75
$ws.Range("M116").Value = -752.75

$ws.Range("H129").Value = 3490.1428
$ws.Range("J129").Value = 4857.75
$ws.Range("L129").Value = 14573.25
$ws.Range("N129").Value = -24573.25

$ws.Range("H131").Value = 51396.95
$ws.Range("I131").Value = 91797.17999999999
$ws.Range("K131").Value = 275391.54
$ws.Range("M131").Value = -270351.54


# ----- Worksheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11602.682
$ws.Range("I122").Value = 12852.833
$ws.Range("K122").Value = 38558.499
$ws.Range("M122").Value = -36108.499

$ws.Range("H126").Value = 3305.8948
$ws.Range("I126").Value = 2799.1428
$ws.Range("K126").Value = 8397.428400000001
$ws.Range("M126").Value = -5927.428400000001

$ws.Range("H132").Value = 6671.3687
$ws.Range("J132").Value = 7785.3335
$ws.Range("L132").Value = 23356.0005
$ws.Range("N132").Value = -28416.0005

$ws.Range("H141").Value = 129966
$ws.Range("J141").Value = 129966
$ws.Range("L141").Value = 129966
$ws.Range("N141").Value = -140326


# ----- Worksheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6174854
$ws.Range("I40").Value = 2266
$ws.Range("J40").Value = 55555556
$ws.Range("K40").Value = 2266
$ws.Range("L40").Value = 55555556
$ws.Range("M40").Value = -2130
$ws.Range("N40").Value = -55555828

$ws.Range("H122").Value = 6264613
$ws.Range("I122").Value = 19903.285
$ws.Range("K122").Value = 59709.855
$ws.Range("M122").Value = -57259.855

$ws.Range("H132").Value = 11815.518
$ws.Range("I132").Value = 17986.295
$ws.Range("K132").Value = 53958.88499999999
$ws.Range("M132").Value = -51428.88499999999

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 115214.5
$ws.Range("J135").Value = 115214.5
$ws.Range("L135").Value = 115214.5
$ws.Range("N135").Value = -125354.5


# ----- Worksheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1740689.6
$ws.Range("I132").Value = 1288.7858
$ws.Range("J132").Value = 3954472.5
$ws.Range("K132").Value = 3866.3574
$ws.Range("L132").Value = 11863417.5
$ws.Range("M132").Value = -1336.3574
$ws.Range("N132").Value = -11868477.5

$ws.Range("H135").Value = 83993.5
$ws.Range("J135").Value = 83993.5
$ws.Range("L135").Value = 83993.5
$ws.Range("N135").Value = -94133.5

$ws.Range("H137").Value = 70000
$ws.Range("J137").Value = 70000
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200

$ws.Range("H140").Value = 84299.336
$ws.Range("J140").Value = 84299.336
$ws.Range("L140").Value = 84299.336
$ws.Range("N140").Value = -94659.336

$ws.Range("H141").Value = 76803.5
$ws.Range("J141").Value = 76803.5
$ws.Range("L141").Value = 76803.5
$ws.Range("N141").Value = -87163.5

